$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.024701
$ws.Range("N2").Value = 9.074103000000001
$ws.Range("O2").Value = 0.1596375877334842
$ws.Range("P2").Value = 0.1596375877334843
$ws.Range("Q2").Value = 0.41671305677
$ws.Range("R2").Value = 3.750417510930001
$ws.Range("S2").Value = 0.0777271876132412
$ws.Range("T2").Value = 0.07772718761324121

$ws.Range("O3").Value = 0.6072559333217162
$ws.Range("P3").Value = 0.6072559333217163
$ws.Range("S3").Value = 0.2956715678850776
$ws.Range("T3").Value = 0.2956715678850776

$ws.Range("M4").Value = 4.368554666666666
$ws.Range("N4").Value = 13.105664
$ws.Range("O4").Value = 0.2305634602787257
$ws.Range("P4").Value = 0.2305634602787257
$ws.Range("Q4").Value = 0.6018557764266667
$ws.Range("R4").Value = 5.41670198784
$ws.Range("S4").Value = 0.1122608377405569
$ws.Range("T4").Value = 0.1122608377405569

$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.04818333333333333
$ws.Range("N5").Value = 0.14455
$ws.Range("O5").Value = 0.002543018666073676
$ws.Range("P5").Value = 0.002543018666073677
$ws.Range("Q5").Value = 0.006638217833333334
$ws.Range("R5").Value = 0.05974396050000001
$ws.Range("S5").Value = 0.001238190151631959
$ws.Range("T5").Value = 0.001238190151631959

$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.024701
$ws.Range("N6").Value = 9.074103000000001
$ws.Range("O6").Value = 0.1596375877334842
$ws.Range("P6").Value = 0.1596375877334843
$ws.Range("Q6").Value = 0.4391402064513333
$ws.Range("R6").Value = 3.952261858062001
$ws.Range("S6").Value = 0.08191040012024305
$ws.Range("T6").Value = 0.08191040012024306

$ws.Range("O7").Value = 0.6072559333217162
$ws.Range("P7").Value = 0.6072559333217163
$ws.Range("S7").Value = 0.3115843654366386
$ws.Range("T7").Value = 0.3115843654366386

$ws.Range("M8").Value = 4.368554666666666
$ws.Range("N8").Value = 13.105664
$ws.Range("O8").Value = 0.2305634602787257
$ws.Range("P8").Value = 0.2305634602787257
$ws.Range("Q8").Value = 0.634247153095111
$ws.Range("R8").Value = 5.708224377855999
$ws.Range("S8").Value = 0.1183026225381688
$ws.Range("T8").Value = 0.1183026225381688

$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.04818333333333333
$ws.Range("N9").Value = 0.14455
$ws.Range("O9").Value = 0.002543018666073676
$ws.Range("P9").Value = 0.002543018666073677
$ws.Range("Q9").Value = 0.006995481188888888
$ws.Range("R9").Value = 0.06295933070000001
$ws.Range("S9").Value = 0.001304828514441717
$ws.Range("T9").Value = 0.001304828514441718
